# "new examples of students"
# Replace the single data row (row 2) with a new student's scores.
# The new example only has 50 subject columns (A:AX) instead of the
# previous 77 (A:BY), so the old row is deleted first (this also drops
# the old row-level/cell-level formatting back to the sheet defaults,
# matching the author's re-export) and then repopulated with the new
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep the workbook's date system explicit (matches the re-saved file's
# workbookPr, even though it is expressed with a different attribute name
# after the round trip).
$wb.Date1904 = $false

# Drop the old row 2 entirely so it comes back with plain/default
# formatting (no explicit style/height) instead of inheriting the old
# header-like formatting.
$ws.Rows.Item(2).Delete()

# New student record (id 119919) across 50 subject columns.
$ws.Range("A2").Value = 119919
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 4
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 4
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 4
$ws.Range("X2").Value = 5
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 2
$ws.Range("AA2").Value = 3
$ws.Range("AB2").Value = 4
$ws.Range("AC2").Value = 3
$ws.Range("AD2").Value = 3
$ws.Range("AE2").Value = 4
$ws.Range("AF2").Value = 3
$ws.Range("AG2").Value = 3
$ws.Range("AH2").Value = 4
$ws.Range("AI2").Value = 4
$ws.Range("AJ2").Value = 4
$ws.Range("AK2").Value = 3
$ws.Range("AL2").Value = 3
$ws.Range("AM2").Value = 4
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 4
$ws.Range("AP2").Value = 3
$ws.Range("AQ2").Value = 4
$ws.Range("AR2").Value = 3
$ws.Range("AS2").Value = 3
$ws.Range("AT2").Value = 4
$ws.Range("AU2").Value = 3
$ws.Range("AV2").Value = 4
$ws.Range("AW2").Value = 4
$ws.Range("AX2").Value = 5

# Leave the sheet selection where the author left it: the now-empty
# columns immediately to the right of the new data (AY2:BY2).
$ws.Range("AY2:BY2").Select()
